$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "ARI" row (4) and "LAGO" row (5): LAGO now appears above ARI.
# Using Cut (via a scratch row far outside the used data) preserves the
# original text-typed cells (account numbers with leading zeros) instead
# of retyping them, which would turn them into numbers.
$scratch = $ws.Range("A500:C500")
$ari = $ws.Range("A4:C4")
$ari.Cut($scratch)
$lago = $ws.Range("A5:C5")
$lago.Cut($ari)
$scratch.Cut($lago)

# ARI's balance changed from 40000 to 28000.
$ws.Range("C5").Value = 28000

# Remove the rows for PEDRO (24220.17) and HEMAT (9927.82), which no
# longer appear in the sheet. Delete bottom-up so row numbers of the
# rows still to be removed don't shift.
$ws.Rows(8).Delete()
$ws.Rows(6).Delete()
